# Daily attendance processing - reorders the names listed in the
# "Recorded By" column (column G) of the "Session Analysis Results" sheet.
#
# Observed rule (derived from the source data):
#   - "<name>, System"                 -> "System, <name>"
#   - "<name>, admin@admin.com"        -> "admin@admin.com, <name>"
#   - "backup@backdoor.com, system, System" -> "system, backup@backdoor.com, System"
#   (i.e. swap the first two comma-separated entries when the first entry is
#    "dnasr281@gmail.com", or — for the 3-entry case — when the first two
#    entries are "backup@backdoor.com" / "system"; all other values are left
#    untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "
    $newVal = $val

    if ($parts.Count -eq 2 -and $parts[0] -eq "dnasr281@gmail.com") {
        $newVal = $parts[1] + ", " + $parts[0]
    } elseif ($parts.Count -eq 3 -and $parts[0] -eq "backup@backdoor.com" -and $parts[1] -eq "system") {
        $newVal = $parts[1] + ", " + $parts[0] + ", " + $parts[2]
    }

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
